# Clear the "TRUE" values out of column B (Is_New) for the rows that no
# longer have a value for that field, leaving those cells blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,7,8,9,10,18,19,25,26,27,28,41,42,49,50,53,54,59,60,73,74,77,78,81,82,83,87,88,89,90,91,92,93)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = ""
}
